$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 12:35"

# --- Row 4: Estados Unidos - refreshed stats ---
$ws.Range("B4").Value = 1593296
$ws.Range("C4").Value = 573
$ws.Range("D4").Value = 370864
$ws.Range("E4").Value = 1127484
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = 94948

# --- Rows 39/40: Rumania overtakes Kuwait in ranking, swapping rows ---
$ws.Range("A39").Value = "Rumania"
$ws.Range("B39").Value = 17585
$ws.Range("C39").Value = 198
$ws.Range("D39").Value = 10581
$ws.Range("E39").Value = 5853
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 1151

$ws.Range("A40").Value = "Kuwait"
$ws.Range("B40").Value = 17568
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 4885
$ws.Range("E40").Value = 12559
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 124

# --- Row 58: Marruecos - refreshed stats ---
$ws.Range("B58").Value = 7185
$ws.Range("C58").Value = 52
$ws.Range("D58").Value = 4212
$ws.Range("E58").Value = 2777
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 196

# --- Row 59: Australia - refreshed stats ---
$ws.Range("D59").Value = 6472
$ws.Range("E59").Value = 509

# --- Row 63: Finlandia - refreshed stats ---
$ws.Range("E63").Value = 1387
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 306

# --- Row 93: Lituania - refreshed stats ---
$ws.Range("B93").Value = 1593
$ws.Range("C93").Value = 16
$ws.Range("E93").Value = 483

# --- Row 103: Hong Kong - refreshed stats ---
$ws.Range("B103").Value = 1064
$ws.Range("C103").Value = 8
$ws.Range("D103").Value = 1029
$ws.Range("E103").Value = 31

# --- Rows 136/137: Etiopia overtakes Estado de Palestina in ranking, swapping rows ---
$ws.Range("A136").Value = "Etiopia"
$ws.Range("C136").Value = 9
$ws.Range("D136").Value = 123
$ws.Range("E136").Value = 270
$ws.Range("H136").Value = 5

$ws.Range("A137").Value = "Estado de Palestina"
$ws.Range("B137").Value = 398
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 346
$ws.Range("E137").Value = 50
$ws.Range("H137").Value = 2

# --- Rows 178/179: Angola overtakes Siria in ranking, swapping rows ---
$ws.Range("A178").Value = "Angola"
$ws.Range("C178").Value = 6
$ws.Range("D178").Value = 17
$ws.Range("E178").Value = 38

$ws.Range("A179").Value = "Siria"
$ws.Range("B179").Value = 58
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 36
$ws.Range("E179").Value = 19
